# Natmi following Dr Hou advice
#
# The LR-pair sheet (Fgf17 -> Fgfr1) is recomputed: the 3 existing rows
# (sender/receiver cluster combinations) get new edge statistics, and the
# "reverse direction" combinations (3 more rows, sender cluster "sCs") are
# appended below, growing the used range from A1:T4 to A1:T7.
#
# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
#          D Target cluster  | E..T numeric edge/expression statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A..D (text) and E..T (numeric) values for data rows 2-7, in sheet order.
$rows = @(
    @{ Row = 2;  A = "FAPs"; B = "Fgf17"; C = "Fgfr1"; D = "ECs";
       Nums = @(1, 0.3333333333333333, 0.05925133333333333, 0.177754,
                0.129311749608618, 0.129311749608618, 3, 1,
                4.675378666666666, 14.026136, 0.03681964474327726,
                0.03681964474327726, 0.2770224198382222, 2.493201778544,
                0.004761212681720937, 0.004761212681720937) },
    @{ Row = 3;  A = "FAPs"; B = "Fgf17"; C = "Fgfr1"; D = "FAPs";
       Nums = @(1, 0.3333333333333333, 0.05925133333333333, 0.177754,
                0.129311749608618, 0.129311749608618, 3, 1,
                82.95722966666666, 248.871689, 0.653306596744776,
                0.653306596744776, 4.915326467389555, 44.237938206506,
                0.08448021905591882, 0.08448021905591882) },
    @{ Row = 4;  A = "FAPs"; B = "Fgf17"; C = "Fgfr1"; D = "sCs";
       Nums = @(1, 0.3333333333333333, 0.05925133333333333, 0.177754,
                0.129311749608618, 0.129311749608618, 3, 1,
                39.34793966666667, 118.043819, 0.3098737585119468,
                0.3098737585119468, 2.331417889169555, 20.982761002526,
                0.04007031787097821, 0.04007031787097821) },
    @{ Row = 5;  A = "sCs";  B = "Fgf17"; C = "Fgfr1"; D = "ECs";
       Nums = @(3, 1, 0.398954, 1.196862,
                0.8706882503913821, 0.8706882503913821, 3, 1,
                4.675378666666666, 14.026136, 0.03681964474327726,
                0.03681964474327726, 1.865261020581333, 16.787349185232,
                0.03205843206155633, 0.03205843206155633) },
    @{ Row = 6;  A = "sCs";  B = "Fgf17"; C = "Fgfr1"; D = "FAPs";
       Nums = @(3, 1, 0.398954, 1.196862,
                0.8706882503913821, 0.8706882503913821, 3, 1,
                82.95722966666666, 248.871689, 0.653306596744776,
                0.653306596744776, 33.09611860443533, 297.8650674399181,
                0.5688263776888572, 0.5688263776888572) },
    @{ Row = 7;  A = "sCs";  B = "Fgf17"; C = "Fgfr1"; D = "sCs";
       Nums = @(3, 1, 0.398954, 1.196862,
                0.8706882503913821, 0.8706882503913821, 3, 1,
                39.34793966666667, 118.043819, 0.3098737585119468,
                0.3098737585119468, 15.69801792177534, 141.282161295978,
                0.2698034406409686, 0.2698034406409686) }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    # Columns E (5) through T (20) hold the numeric statistics.
    for ($i = 0; $i -lt $r.Nums.Length; $i++) {
        $ws.Cells.Item($rowNum, 5 + $i).Value = $r.Nums[$i]
    }
}
